# Fruta / hortaliza, semanal
# Weekly data refresh: the data rows (2-9) get updated with new values
# for Fecha (D), Calidad (L), Volumen (M), Precio mínimo (N), Precio máximo (O),
# Precio promedio ponderado (P), Unidad de comercialización (Q), Precio $/Kg (S)
# and Kg / unidad (T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44650
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 31000
$ws.Range("O2").Value = 32000
$ws.Range("P2").Value = 31500
$ws.Range("S2").Value = 1575

# Row 3
$ws.Range("D3").Value = 44650
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 250

# Row 4
$ws.Range("D4").Value = 44664
$ws.Range("M4").Value = 150
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("S4").Value = 1639
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44636
$ws.Range("L5").Value = "Primera"
$ws.Range("N5").Value = 29000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 29500
$ws.Range("S5").Value = 1475

# Row 6
$ws.Range("D6").Value = 44679
$ws.Range("M6").Value = 200
$ws.Range("Q6").Value = "$/caja 20 kilos"
$ws.Range("S6").Value = 1475
$ws.Range("T6").Value = 20

# Row 7
$ws.Range("D7").Value = 44679
$ws.Range("L7").Value = "Tercera"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 24500
$ws.Range("S7").Value = 1225

# Row 8
$ws.Range("D8").Value = 44643
$ws.Range("N8").Value = 28000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 29000
$ws.Range("S8").Value = 1450

# Row 9
$ws.Range("D9").Value = 44671
$ws.Range("M9").Value = 200
